$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MarksAwarded (column F) values for rows 6-46
$ws.Range("F6").Value = 7
$ws.Range("F7").Value = 5
$ws.Range("F8").Value = 8.5
$ws.Range("F9").Value = 8
$ws.Range("F10").Value = 7
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = 9.5
$ws.Range("F13").Value = 2
$ws.Range("F14").Value = 8.5
$ws.Range("F15").Value = 6.5
$ws.Range("F16").Value = 4.5
$ws.Range("F17").Value = 9.5
$ws.Range("F18").Value = 7.5
$ws.Range("F19").Value = 7
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 2
$ws.Range("F22").Value = 8.5
$ws.Range("F23").Value = 5.5
$ws.Range("F24").Value = 4.5
$ws.Range("F25").Value = 9
$ws.Range("F26").Value = 9
$ws.Range("F27").Value = 6.5
$ws.Range("F28").Value = 7
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 2
$ws.Range("F31").Value = 8.5
$ws.Range("F32").Value = 7
$ws.Range("F33").Value = 6.5
$ws.Range("F34").Value = 9
$ws.Range("F35").Value = 9
$ws.Range("F36").Value = 6.5
$ws.Range("F37").Value = 7
$ws.Range("F38").Value = 2
$ws.Range("F39").Value = 9
$ws.Range("F40").Value = 2
$ws.Range("F41").Value = 9
$ws.Range("F42").Value = 7
$ws.Range("F43").Value = 4.5
$ws.Range("F44").Value = 9
$ws.Range("F45").Value = 7.5
$ws.Range("F46").Value = 7

# Remove the entire Feedback column (G) - no longer used
$ws.Columns.Item(7).Delete()
